$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update loading_percent values for rows 2-25 (case with 380 kV)
$ws.Range("B2").Value2 = 16.28558348298703
$ws.Range("C2").Value2 = 11.21434405186053
$ws.Range("D2").Value2 = 9.880341846853174
$ws.Range("E2").Value2 = 13.95445731972162
$ws.Range("F2").Value2 = 29.4198112808336
$ws.Range("I2").Value2 = 19.27409676538315
$ws.Range("J2").Value2 = 9.787616231031489
$ws.Range("O2").Value2 = 21.72566271065692

$ws.Range("B3").Value2 = 15.53553635930225
$ws.Range("C3").Value2 = 10.57202405173082
$ws.Range("D3").Value2 = 9.811964237362423
$ws.Range("E3").Value2 = 13.88622203461793
$ws.Range("F3").Value2 = 29.50437966604854
$ws.Range("I3").Value2 = 19.44544583590759
$ws.Range("J3").Value2 = 9.795528358734193
$ws.Range("O3").Value2 = 21.84320915618736

$ws.Range("B4").Value2 = 15.05590403970334
$ws.Range("C4").Value2 = 10.15638819930143
$ws.Range("D4").Value2 = 9.771021521453013
$ws.Range("E4").Value2 = 13.8469139762207
$ws.Range("F4").Value2 = 29.566376447605
$ws.Range("I4").Value2 = 19.55627836323962
$ws.Range("J4").Value2 = 9.802050078544363
$ws.Range("O4").Value2 = 21.92237655212741

$ws.Range("B5").Value2 = 14.85587090174414
$ws.Range("C5").Value2 = 9.981770028799989
$ws.Range("D5").Value2 = 9.754611653543574
$ws.Range("E5").Value2 = 13.83155875148783
$ws.Range("F5").Value2 = 29.59416026829762
$ws.Range("I5").Value2 = 19.6028589652539
$ws.Range("J5").Value2 = 9.805126129649066
$ws.Range("O5").Value2 = 21.95638776736272

$ws.Range("B6").Value2 = 14.82238620963235
$ws.Range("C6").Value2 = 9.952461515921398
$ws.Range("D6").Value2 = 9.751903758068801
$ws.Range("E6").Value2 = 13.82904941785813
$ws.Range("F6").Value2 = 29.59892553358095
$ws.Range("I6").Value2 = 19.6106791704042
$ws.Range("J6").Value2 = 9.805662175315286
$ws.Range("O6").Value2 = 21.96214073579848

$ws.Range("B7").Value2 = 15.05322454225485
$ws.Range("C7").Value2 = 10.15405431818304
$ws.Range("D7").Value2 = 9.770799083773788
$ws.Range("E7").Value2 = 13.84670418987347
$ws.Range("F7").Value2 = 29.56674096587949
$ws.Range("I7").Value2 = 19.55690083258491
$ws.Range("J7").Value2 = 9.802089869175546
$ws.Range("O7").Value2 = 21.92282816572728

$ws.Range("B8").Value2 = 16.03105545504938
$ws.Range("C8").Value2 = 10.99733553633764
$ws.Range("D8").Value2 = 9.856556874288616
$ws.Range("E8").Value2 = 13.93040008656536
$ws.Range("F8").Value2 = 29.44687265140214
$ws.Range("I8").Value2 = 19.33201157719906
$ws.Range("J8").Value2 = 9.789999111785466
$ws.Range("O8").Value2 = 21.76473603904745

$ws.Range("B9").Value2 = 17.78882160445054
$ws.Range("C9").Value2 = 12.47904981584828
$ws.Range("D9").Value2 = 10.03241323302141
$ws.Range("E9").Value2 = 14.1144928104993
$ws.Range("F9").Value2 = 29.2922565372358
$ws.Range("I9").Value2 = 18.93556228612077
$ws.Range("J9").Value2 = 9.779483241160925
$ws.Range("O9").Value2 = 21.51059287121541

$ws.Range("B10").Value2 = 18.97385515298161
$ws.Range("C10").Value2 = 13.45949919669939
$ws.Range("D10").Value2 = 10.16550534432551
$ws.Range("E10").Value2 = 14.26109401486251
$ws.Range("F10").Value2 = 29.22836151625173
$ws.Range("I10").Value2 = 18.67136900038356
$ws.Range("J10").Value2 = 9.779787126839659
$ws.Range("O10").Value2 = 21.358459229459

$ws.Range("B11").Value2 = 19.48843426857424
$ws.Range("C11").Value2 = 13.88159880988164
$ws.Range("D11").Value2 = 10.22672274331512
$ws.Range("E11").Value2 = 14.33006039955103
$ws.Range("F11").Value2 = 29.2102032980803
$ws.Range("I11").Value2 = 18.55704456381945
$ws.Range("J11").Value2 = 9.781663632643593
$ws.Range("O11").Value2 = 21.29686976760751

$ws.Range("B12").Value2 = 19.67967655310762
$ws.Range("C12").Value2 = 14.03797336093675
$ws.Range("D12").Value2 = 10.24998601723045
$ws.Range("E12").Value2 = 14.35648653858966
$ws.Range("F12").Value2 = 29.20490337592393
$ws.Range("I12").Value2 = 18.51459456962973
$ws.Range("J12").Value2 = 9.782623474346057
$ws.Range("O12").Value2 = 21.27465169662614

$ws.Range("B13").Value2 = 19.63865133637576
$ws.Range("C13").Value2 = 14.00444977722247
$ws.Range("D13").Value2 = 10.24497247008266
$ws.Range("E13").Value2 = 14.35078168794684
$ws.Range("F13").Value2 = 29.20597459952264
$ws.Range("I13").Value2 = 18.52369948099735
$ws.Range("J13").Value2 = 9.782405682567909
$ws.Range("O13").Value2 = 21.27938748656659

$ws.Range("B14").Value2 = 19.50424087397989
$ws.Range("C14").Value2 = 13.89453341762394
$ws.Range("D14").Value2 = 10.22863507945392
$ws.Range("E14").Value2 = 14.3322283808575
$ws.Range("F14").Value2 = 29.20973564183886
$ws.Range("I14").Value2 = 18.55353529936007
$ws.Range("J14").Value2 = 9.781737608613534
$ws.Range("O14").Value2 = 21.29501968497905

$ws.Range("B15").Value2 = 19.42143692880068
$ws.Range("C15").Value2 = 13.82675450582894
$ws.Range("D15").Value2 = 10.21863812618563
$ws.Range("E15").Value2 = 14.3209038009495
$ws.Range("F15").Value2 = 29.21224486089149
$ws.Range("I15").Value2 = 18.5719202729268
$ws.Range("J15").Value2 = 9.781360830297062
$ws.Range("O15").Value2 = 21.30473897039237

$ws.Range("B16").Value2 = 18.93972473454815
$ws.Range("C16").Value2 = 13.43143034801721
$ws.Range("D16").Value2 = 10.16151688917453
$ws.Range("E16").Value2 = 14.25663124464894
$ws.Range("F16").Value2 = 29.2297683789038
$ws.Range("I16").Value2 = 18.67895823294981
$ws.Range("J16").Value2 = 9.779699413437424
$ws.Range("O16").Value2 = 21.36263830853123

$ws.Range("B17").Value2 = 18.6378623396121
$ws.Range("C17").Value2 = 13.18276571476319
$ws.Range("D17").Value2 = 10.12663646429649
$ws.Range("E17").Value2 = 14.21777309098703
$ws.Range("F17").Value2 = 29.2433181435504
$ws.Range("I17").Value2 = 18.74612281110439
$ws.Range("J17").Value2 = 9.779124936609159
$ws.Range("O17").Value2 = 21.40011515194632

$ws.Range("B18").Value2 = 18.46193765244658
$ws.Range("C18").Value2 = 13.0374907647873
$ws.Range("D18").Value2 = 10.10663844297294
$ws.Range("E18").Value2 = 14.19563828240582
$ws.Range("F18").Value2 = 29.25213794290227
$ws.Range("I18").Value2 = 18.78530552638706
$ws.Range("J18").Value2 = 9.778958177653395
$ws.Range("O18").Value2 = 21.42238723592428

$ws.Range("B19").Value2 = 18.40198047193272
$ws.Range("C19").Value2 = 12.98791773696115
$ws.Range("D19").Value2 = 10.09987895845273
$ws.Range("E19").Value2 = 14.1881813262598
$ws.Range("F19").Value2 = 29.2553001941858
$ws.Range("I19").Value2 = 18.79866683862821
$ws.Range("J19").Value2 = 9.77892984500836
$ws.Range("O19").Value2 = 21.4300509391686

$ws.Range("B20").Value2 = 18.67023503511527
$ws.Range("C20").Value2 = 13.20946953730348
$ws.Range("D20").Value2 = 10.13034300269978
$ws.Range("E20").Value2 = 14.22188744125127
$ws.Range("F20").Value2 = 29.24176946662704
$ws.Range("I20").Value2 = 18.73891595877167
$ws.Range("J20").Value2 = 9.77916915666585
$ws.Range("O20").Value2 = 21.39605146651627

$ws.Range("B21").Value2 = 19.54381934076162
$ws.Range("C21").Value2 = 13.92691274413867
$ws.Range("D21").Value2 = 10.2334316748966
$ws.Range("E21").Value2 = 14.33766966249884
$ws.Range("F21").Value2 = 29.2085881005925
$ws.Range("I21").Value2 = 18.54474894066046
$ws.Range("J21").Value2 = 9.781927079896569
$ws.Range("O21").Value2 = 21.29039808020709

$ws.Range("B22").Value2 = 20.09363790688439
$ws.Range("C22").Value2 = 14.37559686338581
$ws.Range("D22").Value2 = 10.30127445512182
$ws.Range("E22").Value2 = 14.41513894711355
$ws.Range("F22").Value2 = 29.19609205512402
$ws.Range("I22").Value2 = 18.42275849307607
$ws.Range("J22").Value2 = 9.785181935131316
$ws.Range("O22").Value2 = 21.22778965276528

$ws.Range("B23").Value2 = 19.80214895848601
$ws.Range("C23").Value2 = 14.13798159701819
$ws.Range("D23").Value2 = 10.26502767497717
$ws.Range("E23").Value2 = 14.37363338836145
$ws.Range("F23").Value2 = 29.20191837587056
$ws.Range("I23").Value2 = 18.48741796266713
$ws.Range("J23").Value2 = 9.783312133247582
$ws.Range("O23").Value2 = 21.26061248202437

$ws.Range("B24").Value2 = 18.65560674522317
$ws.Range("C24").Value2 = 13.19740394574128
$ws.Range("D24").Value2 = 10.12866710360462
$ws.Range("E24").Value2 = 14.22002670283566
$ws.Range("F24").Value2 = 29.24246641606285
$ws.Range("I24").Value2 = 18.74217240451532
$ws.Range("J24").Value2 = 9.779148655389941
$ws.Range("O24").Value2 = 21.39788639743488

$ws.Range("B25").Value2 = 17.33143084942764
$ws.Range("C25").Value2 = 12.09701482776449
$ws.Range("D25").Value2 = 9.98409755898329
$ws.Range("E25").Value2 = 14.06263536441656
$ws.Range("F25").Value2 = 29.32539944408652
$ws.Range("I25").Value2 = 19.03805146126546
$ws.Range("J25").Value2 = 9.780916447615102
$ws.Range("O25").Value2 = 21.57330925135392
